# "Generate Report for Handback" - mark the two reports' entries as handed
# back (in sync with en-US), stamp the handback datetime, and link the
# "Latest Target File" column to the (localized) a.md source file.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5237317f0d835b0a9dd99677f083d0a40788af52/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status summary columns (E, F) for both
# rows now reflect the handed-back status.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Widen the zh-cn / de-de columns to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Per-locale detail sheets (zh-cn, de-de): Status column, the newly
# filled-in "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns, and a hyperlink on the target file.
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Handback = "2016-08-28 00:34:59" },
    @{ Sheet = "de-de"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Handback = "2016-08-28 00:35:11" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Status column (C) for both data rows.
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File (I): link to the source a.md that was handed back.
    $ws.Hyperlinks.Add($ws.Range("I2"), $aMdUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $aMdUrl, "", "", "a.md")

    # Latest Handback File (J).
    $ws.Range("J2").Value = $locale.Xlf
    $ws.Range("J3").Value = $locale.Xlf

    # Latest Handback DateTime (K).
    $ws.Range("K2").Value = $locale.Handback
    $ws.Range("K3").Value = $locale.Handback

    # Widen Status (C) and Latest Handback File (J) columns for the
    # longer text they now hold.
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}
